$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D): numeric-looking text values that must stay stored as text ---
# Force text format before assignment so Excel does not coerce the numeric-looking
# strings into actual numbers, then restore the default "Normal" style so no stray
# number-format styling is left behind on the cell.
$priceCells = @("D2", "D3", "D4", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D25", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '246.21'
$ws.Range("D3").Value = '24.14'
$ws.Range("D4").Value = '5.281'
$ws.Range("D6").Value = '6.501'
$ws.Range("D7").Value = '3.133'
$ws.Range("D8").Value = '0.8171'
$ws.Range("D9").Value = '0.8555'
$ws.Range("D10").Value = '0.1358'
$ws.Range("D11").Value = '0.06933'
$ws.Range("D12").Value = '0.03129'
$ws.Range("D13").Value = '0.02874'
$ws.Range("D14").Value = '0.09406'
$ws.Range("D15").Value = '3.748'
$ws.Range("D16").Value = '0.001511'
$ws.Range("D17").Value = '0.04670'
$ws.Range("D18").Value = '0.0005987'
$ws.Range("D19").Value = '0.006271'
$ws.Range("D20").Value = '0.001236'
$ws.Range("D21").Value = '0.004629'
$ws.Range("D22").Value = '0.00006894'
$ws.Range("D23").Value = '3.498'
$ws.Range("D25").Value = '0.3192'
$ws.Range("D40").Value = '0.03662'
$ws.Range("D41").Value = '0.1057'
$ws.Range("D42").Value = '0.002748'
$ws.Range("D43").Value = '0.003016'
$ws.Range("D44").Value = '0.007484'
$ws.Range("D45").Value = '0.00005262'
$ws.Range("D46").Value = '0.00000000750'
$ws.Range("D47").Value = '0.3698'
$ws.Range("D48").Value = '0.002237'
$ws.Range("D49").Value = '0.00002099'
$ws.Range("D50").Value = '0.0001999'

foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}

# --- Coin / Link / Rank-label columns (B, C, E): plain text, no special handling needed ---
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("E16").Value = '15BitForexTokenBF'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("E18").Value = '17OneONE'
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("E41").Value = '40BKEXTokenBKK'
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("E42").Value = '41CEJICEJIBestin24h'
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("E43").Value = '42KickTokenKICKWorstin24h'
$ws.Range("E47").Value = '46CoinbaseStockTokenCOIN'
